$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, preserving the original (default) cell style
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.091.85'
Set-TextValue 'E2' '  -1.79%  '
Set-TextValue 'D3' '1.834.62'
Set-TextValue 'E3' '  -1.37%  '
Set-TextValue 'D4' '0.9993'
Set-TextValue 'E4' '  -0.05%  '
Set-TextValue 'D5' '239.73'
Set-TextValue 'E5' '  -2.23%  '
Set-TextValue 'D6' '0.6792'
Set-TextValue 'E6' '  -2.63%  '
Set-TextValue 'D7' '1.0000'
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'D8' '0.2979'
Set-TextValue 'E8' '  -2.76%  '
Set-TextValue 'D9' '0.07447'
Set-TextValue 'E9' '  -3.51%  '
Set-TextValue 'D10' '23.14'
Set-TextValue 'E10' '  -2.11%  '
Set-TextValue 'D11' '0.07648'
Set-TextValue 'E11' '  -1.38%  '
Set-TextValue 'D12' '1.838.61'
Set-TextValue 'E12' '  -1.17%  '
Set-TextValue 'D13' '5.020'
Set-TextValue 'E13' '  -2.74%  '
Set-TextValue 'D14' '0.6754'
Set-TextValue 'E14' '  -2.43%  '
Set-TextValue 'D15' '86.58'
Set-TextValue 'E15' '  -6.14%  '
Set-TextValue 'D16' '6.159'
Set-TextValue 'E16' '  -6.02%  '
Set-TextValue 'D17' '29.099.02'
Set-TextValue 'D18' '0.000008261'
Set-TextValue 'E18' '  -0.98%  '
Set-TextValue 'D19' '2.078.37'
Set-TextValue 'E19' '  -1.25%  '
Set-TextValue 'D20' '227.85'
Set-TextValue 'E20' '  -5.72%  '
Set-TextValue 'D21' '12.45'
Set-TextValue 'E21' '  -2.34%  '
Set-TextValue 'D22' '0.9995'
Set-TextValue 'E22' '  -0.05%  '
Set-TextValue 'D23' '7.327'
Set-TextValue 'E23' '  -3.70%  '
Set-TextValue 'E24' '  -0.04%  '
Set-TextValue 'D25' '160.83'
Set-TextValue 'E25' '  +0.73%  '
Set-TextValue 'D26' '0.1438'
Set-TextValue 'E26' '  -4.45%  '
Set-TextValue 'E27' '  -2.47%  '
Set-TextValue 'D28' '17.99'
Set-TextValue 'E28' '  -1.64%  '
Set-TextValue 'D29' '1.498'
Set-TextValue 'E29' '  -2.57%  '
Set-TextValue 'D30' '4.243'
Set-TextValue 'E30' '  -0.15%  '
Set-TextValue 'E31' '  -1.41%  '
Set-TextValue 'D32' '1.195'
Set-TextValue 'E32' '  +0.17%  '
Set-TextValue 'D33' '0.05372'
Set-TextValue 'E33' '  +5.33%  '
Set-TextValue 'D34' '0.7529'
Set-TextValue 'E34' '  -3.32%  '
Set-TextValue 'D35' '1.847'
Set-TextValue 'E35' '  -2.69%  '
Set-TextValue 'E36' '  -2.18%  '
Set-TextValue 'D37' '2.682'
Set-TextValue 'E37' '  -0.11%  '
Set-TextValue 'D38' '1.308.19'
Set-TextValue 'E38' '  -1.35%  '
Set-TextValue 'D39' '0.01813'
Set-TextValue 'E39' '  -3.25%  '
Set-TextValue 'D40' '2.714'
Set-TextValue 'E40' '  -0.73%  '
Set-TextValue 'D41' '0.9333'
Set-TextValue 'E41' '  -2.64%  '
Set-TextValue 'D42' '6.051'
Set-TextValue 'E42' '  +3.61%  '
Set-TextValue 'D43' '0.08587'
Set-TextValue 'E43' '  +35.06%  '
Set-TextValue 'D44' '104.87'
Set-TextValue 'E44' '  -1.30%  '
Set-TextValue 'D45' '0.9989'
Set-TextValue 'E45' '  -0.09%  '
Set-TextValue 'D46' '1.979.96'
Set-TextValue 'E46' '  -1.25%  '
Set-TextValue 'D47' '0.5178'
Set-TextValue 'E47' '  -0.66%  '
Set-TextValue 'D48' '1.767'
Set-TextValue 'E48' '  -0.91%  '
Set-TextValue 'B49' 'Aave'
Set-TextValue 'C49' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D49' '63.72'
Set-TextValue 'E49' '  -1.25%  '
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '9.379'
Set-TextValue 'E50' '  -3.99%  '
Set-TextValue 'B51' 'Cronos'
Set-TextValue 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D51' '0.05928'
Set-TextValue 'E51' '  +0.21%  '
